$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "city" column before the old column B (Temparature),
# shifting Temparature/Windspeed/Event one column to the right.
$ws.Columns("B").Insert()

# New header
$ws.Range("B1").Value = "city"

# City values for the existing rows (HYD / RCB groups)
$ws.Range("B2").Value = "HYD"
$ws.Range("B3").Value = "HYD"
$ws.Range("B4").Value = "HYD"
$ws.Range("B5").Value = "RCB"
$ws.Range("B6").Value = "RCB"
$ws.Range("B7").Value = "RCB"

# Give the new column the same width Excel shows for column A
$ws.Columns("B").ColumnWidth = 9.5

# New rows of weather data (CSK), reusing the date-formatted style from
# the existing date column so the new dates render the same way
$ws.Range("A2").Copy()
$ws.Range("A8:B10").PasteSpecial(-4122)

$ws.Range("A8").Value = 42917
$ws.Range("B8").Value = "CSK"
$ws.Range("C8").Value = 29
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = "sunny"

$ws.Range("A9").Value = 42948
$ws.Range("B9").Value = "CSK"
$ws.Range("C9").Value = 34
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = "snow"

$ws.Range("A10").Value = 42979
$ws.Range("B10").Value = "CSK"
$ws.Range("C10").Value = 35
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = "rainy"

# Bold the header row
$ws.Range("A1:E1").Font.Bold = $true

# Update the active selection to match where Excel left it
$ws.Range("E10").Select()
